# Season update through 1/17:
# The game played on 2024-01-15 (serial 45306) vs PHI — previously listed
# as an upcoming game on the "Next" sheet — has now been played, so:
#   1. Its final box-score line is appended to the "Games" sheet as row 40
#      (game #39).
#   2. It is removed from the "Next" sheet, shifting the remaining
#      upcoming games up by one row.

$wb = $excel.ActiveWorkbook

# ---- 1. Append the completed game to the "Games" sheet ----
$games = $wb.Worksheets.Item("Games")

$newRow = 40

$games.Cells.Item($newRow, 1).Value = 39
$games.Cells.Item($newRow, 2).Value = 45306
$games.Cells.Item($newRow, 2).NumberFormat = "YYYY-MM-DD"
$games.Cells.Item($newRow, 3).Value = -2
$games.Cells.Item($newRow, 4).Value = 115
$games.Cells.Item($newRow, 5).Value = 97.7
$games.Cells.Item($newRow, 6).Value = 0.545
$games.Cells.Item($newRow, 7).Value = 11.3
$games.Cells.Item($newRow, 8).Value = 32.7
$games.Cells.Item($newRow, 9).Value = 0.216
$games.Cells.Item($newRow, 10).Value = 117.7
$games.Cells.Item($newRow, 11).Value = "PHI"
$games.Cells.Item($newRow, 12).Value = 124
$games.Cells.Item($newRow, 13).Value = 0.572
$games.Cells.Item($newRow, 14).Value = 9.300000000000001
$games.Cells.Item($newRow, 15).Value = 21.1
$games.Cells.Item($newRow, 16).Value = 0.349
$games.Cells.Item($newRow, 17).Value = 127
$games.Cells.Item($newRow, 18).Value = 0
$games.Cells.Item($newRow, 19).Value = 0

# ---- 2. Remove the now-played game from the "Next" sheet ----
$next = $wb.Worksheets.Item("Next")
$next.Rows.Item(2).Delete()
